# Add a new conference entry (row 20) to the conferences table:
#   University of Durham - Conference of the UK Algebraic Geometry Network
#   11th - 12th September, 2024 - hyperlinked to the event page.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 20

$ws.Cells.Item($row, 1).Value = -1.57217905973651
$ws.Cells.Item($row, 2).Value = 54.7638082226076
$ws.Cells.Item($row, 3).Value = "University of Durham"
$ws.Cells.Item($row, 4).Value = "Conference of the UK Algebraic Geometry Network"
$ws.Cells.Item($row, 5).Value = "11th - 12th September, 2024"

$linkUrl = "https://www.ukagnetwork.org/upcoming-activities/durham-september-2024"
$linkCell = $ws.Cells.Item($row, 6)
$linkCell.Value = $linkUrl

# Register the hyperlink relationship (this also drags in Excel's built-in
# "Hyperlink" style, which does not match the workbook's existing look).
$ws.Hyperlinks.Add($linkCell, $linkUrl, "", "", $linkUrl) | Out-Null

# The rest of the "Link" column (F2:F19) uses plain direct formatting (blue
# Arial, no underline) rather than the Hyperlink theme style, so re-apply
# that existing look to the new cell by cloning the format from the row
# above, putting the cell style back in line with the rest of the sheet.
$ws.Range("F19").Copy() | Out-Null
$linkCell.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Drop the now-unused built-in "Hyperlink" cell style that Hyperlinks.Add
# registered, so the style table stays the way it was before the edit.
$wb.Styles.Item("Hyperlink").Delete()

# Keep the view roughly where the author left it after adding the row.
$ws.Range("A1").Select() | Out-Null
$ws.Range("A20:B20").Select() | Out-Null
